$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4543
$ws.Range("F5").Value = 3656
$ws.Range("F6").Value = 1056
$ws.Range("F7").Value = 168
$ws.Range("F9").Value = 361
$ws.Range("F10").Value = 359
$ws.Range("F11").Value = 2522
$ws.Range("F17").Value = 553
$ws.Range("F19").Value = 62
$ws.Range("F20").Value = 10445
$ws.Range("F21").Value = 6066
$ws.Range("F28").Value = 840
$ws.Range("F30").Value = 175
$ws.Range("F35").Value = 481
$ws.Range("F37").Value = 269
$ws.Range("F39").Value = 245
$ws.Range("F40").Value = 4852
$ws.Range("F41").Value = 26
$ws.Range("F42").Value = 1140
$ws.Range("F44").Value = 183
$ws.Range("F46").Value = 489

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 30
$ws.Range("F15").Value = 3580

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8811
$ws.Range("F3").Value = 442
$ws.Range("F4").Value = 1642

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 442
$ws.Range("F4").Value = 4543
$ws.Range("F7").Value = 3656
$ws.Range("F8").Value = 168
$ws.Range("F10").Value = 2522
$ws.Range("F11").Value = 30
$ws.Range("F17").Value = 553
$ws.Range("F19").Value = 62
$ws.Range("F20").Value = 10446
$ws.Range("F21").Value = 3580
$ws.Range("F29").Value = 840
$ws.Range("F31").Value = 175
$ws.Range("F35").Value = 481
$ws.Range("F37").Value = 269
$ws.Range("F40").Value = 245
$ws.Range("F41").Value = 4852
$ws.Range("F42").Value = 26
$ws.Range("F43").Value = 1140
$ws.Range("F45").Value = 489
